# Insert a new data row (weekly Cilantro price record) before the current
# row 104, shifting all subsequent rows (104-186) down by one (105-187).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(104).Insert()

# Populate the newly inserted row 104 with the new record.
$ws.Range("A104").Value = 3
$ws.Range("B104").Value = "Femacal de La Calera"
$ws.Range("C104").Value = "Coquimbo"
$ws.Range("D104").Value = (Get-Date -Year 2021 -Month 9 -Day 8).Date
$ws.Range("E104").Value = 5
$ws.Range("F104").Value = 100112040
$ws.Range("G104").Value = "Cilantro"
$ws.Range("H104").Value = "Sin especificar"
$ws.Range("I104").Value = "Primera"
$ws.Range("J104").Value = 310
$ws.Range("K104").Value = 2500
$ws.Range("L104").Value = 3000
$ws.Range("M104").Value = 2758
$ws.Range("N104").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O104").Value = "Provincia de Quillota"
$ws.Range("P104").Value = 919
$ws.Range("Q104").Value = 3
$ws.Range("R104").Value = "Hortaliza"
